# edit.ps1
# Applies the "Save system + minor fixes" commit to To-Do-List.docx:
#   - Fills in the text of the existing trailing (empty) list item with
#     "Add more sounds effects (button presses)".
#   - Appends nine further ListParagraph bullet items describing further
#     to-do work, preserving the document's hidden "_GoBack" bookmark by
#     re-homing it (via raw OOXML) into the middle of the
#     "More visual upgrades" run split, exactly as in the target revision.
#   - The last four new items ("Save system", "End-game / monthly grade",
#     "Tutorial", "Define smog mechanics") use bold (and, for "Save
#     system", strikethrough) run/paragraph-mark formatting.

$d = $word.ActiveDocument

# 1) The document currently ends with an empty ListParagraph bullet that
#    only carries the leftover "_GoBack" bookmark. Give it its text.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.Text = "Add more sounds effects (button presses)"

# 2) Setting .Range.Text above left the "_GoBack" bookmark sitting at the
#    end of that paragraph. The target document instead has it further
#    down, embedded inside the "More visual upgrades" paragraph (split
#    across runs, right where the cursor last was). Drop the stale
#    bookmark here; the replacement XML below defines it at the right
#    spot.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 3) Append the remaining nine list paragraphs as literal OOXML so the
#    run splits, preserved/xml:space, inline bookmark placement, and
#    bold/strike formatting all come out byte-for-byte as in the target.
$xmlPayload = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Add </w:t></w:r><w:r><w:t xml:space="preserve">“ground”-box (opposite of skybox) </w:t></w:r><w:r><w:t xml:space="preserve">and </w:t></w:r><w:r><w:t>new grass texture</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Make the upgrade tree more logical, balanced, and scrollable </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Rearrange and expand the city, creating a linear path of unlockable areas (delineated by smog)</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Overhaul the monthly report menu</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>More</w:t></w:r><w:r><w:t xml:space="preserve"> visual</w:t></w:r><w:r><w:t xml:space="preserve"> upgr</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>ades</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/><w:strike/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:strike/></w:rPr><w:t>Save system</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t>End-game / monthly grade</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Tutorial </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Define smog mechanics </w:t></w:r></w:p>
'@

$endOfStory = $d.Content.End
$insertionPoint = $d.Range($endOfStory - 1, $endOfStory - 1)
$insertionPoint.InsertXML($xmlPayload)
